$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 6: Pagos (F) and Inscricoes homologadas (H)
$ws.Range("F6").Value = 46
$ws.Range("H6").Value = 56

# Row 12
$ws.Range("F12").Value = 613
$ws.Range("H12").Value = 699

# Row 25
$ws.Range("F25").Value = 256
$ws.Range("H25").Value = 316

# Row 26
$ws.Range("F26").Value = 184
$ws.Range("H26").Value = 209

# Row 28
$ws.Range("F28").Value = 156
$ws.Range("H28").Value = 208

# Row 33
$ws.Range("F33").Value = 231
$ws.Range("H33").Value = 322

# Row 43
$ws.Range("F43").Value = 108
$ws.Range("H43").Value = 135

# Row 46
$ws.Range("F46").Value = 301
$ws.Range("H46").Value = 365

$wb.Save()
